# Google Ads.xlsx — reformat the single "questions" blob from a
# Python-repr single-line string into a pretty-printed, double-quoted
# JSON-ish multi-line string, and drop the old helper cell (A1 = 0)
# that previously held a bold/bordered "header" style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be:
#   A1 (bold, bordered) = 0
#   A2 (shared string)  = questions = [...] (single line, single quotes)
# Delete row 1 entirely so the old A2 content shifts up to A1, and the
# stale bold/border style (style index 1) is no longer referenced by
# any cell.
$ws.Rows.Item(1).Delete()

# Make sure A1 carries no leftover formatting (default style only).
$ws.Range("A1").ClearFormats()

# Write the reformatted / pretty-printed questions text into A1.
$ws.Range("A1").Value = 'questions = [
    {
        "title": "What metric would you use to determine how often your ad is shown compared to the number of impressions you are eligible to receive?",
        "ques_type": 2,
        "options": [
            "Search impression share",
            "Impressions",
            "Search top impression share",
            "Average position"
        ],
        "score": "Search impression share"
    },
    {
        "title": "What can you do with vCPM bidding?",
        "ques_type": 2,
        "options": [
            "You bid for your ad based on how often it appears in a viewable position on the Google Display Network.",
            "You bid for your ad based on how often it appears on the Google Display Network.",
            "You can maximize the number of impressions you can reach daily.",
            "You can set the max amount you want to pay for each ad click."
        ],
        "score": "You bid for your ad based on how often it appears in a viewable position on the Google Display Network."
    },
    {
        "title": "What is the division of credits for Position-based attribution?",
        "ques_type": 2,
        "options": [
            "40% to the first click, 40% to the last click, 20% to the clicks in-between",
            "25% to the first click, 25% to the last click, 50% to the clicks in-between",
            "25% to the first click, 50% to the last click, 25% to the clicks in-between",
            "20% to the first click, 40% to the last click, 40% to the clicks in-between"
        ],
        "score": "40% to the first click, 40% to the last click, 20% to the clicks in-between"
    },
    {
        "title": "What can you do to optimize a search campaign given the following details?Target CPA bidding is activeThe ad is not receiving any impressionsSearch impression share lost (rank) &gt 0",
        "ques_type": 2,
        "options": [
            "Increase target CPA",
            "Decrease target CPA",
            "Increase CPC bids",
            "Set your bid adjustments to bid higher during the best hours of the day."
        ],
        "score": "Increase target CPA"
    }
]'

# Restore the row to automatic (non-custom) height instead of leaving
# an explicit autofit height behind from the multi-line text entry.
$ws.Rows.Item(1).AutoFit()
